$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    '骚扰',
    '气不忿儿',
    '煽动者',
    '初衷',
    '节约',
    '助力',
    '推动',
    '制止',
    '歉意',
    '周全',
    '完善',
    '精准',
    '搭配'
)

$newDefs = @(
    '1.) vt. harass',
    '1.) v. be jealous, take another''s success badly; 2.) v. be unable to contain one''s anger',
    '1.) n. demagogue',
    '1.) n. original intent/aspiration',
    '1.) v. economize, conserve; 2.) adj. frugal, economic',
    '1.) n. a helping hand, help, assistance',
    '1.) v. push forward, promote, to push [for acceptance], to actuate',
    '1.) v. put a stop to, curb, to check, to limit',
    '1.) n. apology; 2.) regret',
    '1.) adj. thorough, comprehensive',
    '1.) v. to perfect, to improve; 2.) adj. perfect',
    '1.) adj. accurate, precise, exact; 2.) precision, accuracy',
    '1.) v. pair up, match, arrange in pairs, add sth into a group'
)

$newDate = "2020-12-04"

$startRow = 76
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newWords[$i]
    $ws.Cells.Item($row, 2).Value = $newDefs[$i]

    # Column C holds the date as literal text (matches the rest of the
    # column, which stores "yyyy-mm-dd" strings rather than real date
    # serials). Force text via NumberFormat, assign, then restore the
    # default "Normal" style so the cell keeps no explicit format -
    # only the value stays textual.
    $cell = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $newDate
    $cell.Style = "Normal"
}
